$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style of existing header cell (e.g. E1) onto the new headers so they
# keep the bold/centered/bordered look.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Outlier flags (MAD-based) per algorithm/row
$flagRows = @(
    @{ Row = 2;  F = $false; G = $false; H = $false },
    @{ Row = 3;  F = $false; G = $false; H = $false },
    @{ Row = 4;  F = $false; G = $false; H = $false },
    @{ Row = 5;  F = $false; G = $false; H = $false },
    @{ Row = 6;  F = $false; G = $false; H = $false },
    @{ Row = 7;  F = $true;  G = $false; H = $false },
    @{ Row = 8;  F = $false; G = $false; H = $false },
    @{ Row = 9;  F = $false; G = $false; H = $false },
    @{ Row = 10; F = $false; G = $false; H = $false },
    @{ Row = 11; F = $true;  G = $false; H = $false },
    @{ Row = 12; F = $false; G = $false; H = $false },
    @{ Row = 13; F = $true;  G = $false; H = $false },
    @{ Row = 14; F = $false; G = $false; H = $false },
    @{ Row = 15; F = $false; G = $false; H = $false },
    @{ Row = 16; F = $false; G = $false; H = $false },
    @{ Row = 17; F = $false; G = $false; H = $false },
    @{ Row = 18; F = $true;  G = $false; H = $false }
)

foreach ($entry in $flagRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
}
